# Weekly CompStat update: new crime data collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell as genuine TEXT
# (shared-string), not an auto-coerced number, while keeping the cell's
# original "General" style/format (matches cells like C14/A14 etc.).
# 1) Temporarily flip the target to a text format so the literal string
#    sticks instead of being parsed as a number.
# 2) Copy the General-formatted style from a known-good donor cell over it
#    with PasteSpecial(formats only) so the final number format/style index
#    matches the rest of the sheet (content is untouched by this paste).
function Set-TextValue($cellRef, $text, $formatDonorRef) {
    $target = $ws.Range($cellRef)
    $target.NumberFormat = "@"
    $target.Value = $text
    $ws.Range($formatDonorRef).Copy()
    $target.PasteSpecial(-4122)
}

# --- Header text updates -------------------------------------------------
# "Volume 30   Number  41" -> "...42"
$ws.Range("A8").Value = "Volume 30   Number  42"
# "Report Covering the Week  10/9/2023  Through  10/15/2023"
#   -> "...10/16/2023  Through  10/22/2023"
$ws.Range("C9").Value = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# --- Crime Complaints table (rows 14-29) ---------------------------------

# Row 14 - Murder
$ws.Range("N14").Value = -77.777777777777

# Row 15 - Rape
Set-TextValue "C15" "0" "C14"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = -6.25
$ws.Range("N15").Value = -28.571428571428

# Row 16 - Robbery
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 118
$ws.Range("J16").Value = 97
$ws.Range("K16").Value = 21.649484536082
$ws.Range("L16").Value = 42.168674698795
$ws.Range("M16").Value = -17.482517482517
$ws.Range("N16").Value = -80.781758957654

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 178
$ws.Range("J17").Value = 147
$ws.Range("K17").Value = 21.088435374149
$ws.Range("L17").Value = 25.352112676056
$ws.Range("M17").Value = 81.632653061224
$ws.Range("N17").Value = -30.739299610894

# Row 18 - Burglary
$ws.Range("C18").Value = 5
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -45
$ws.Range("I18").Value = 175
$ws.Range("J18").Value = 157
$ws.Range("K18").Value = 11.464968152866
$ws.Range("L18").Value = 52.173913043478
$ws.Range("M18").Value = -24.892703862660
$ws.Range("N18").Value = -86.559139784946

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -38.888888888888
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 67
$ws.Range("H19").Value = -50.746268656716
$ws.Range("I19").Value = 522
$ws.Range("J19").Value = 620
$ws.Range("K19").Value = -15.806451612903
$ws.Range("L19").Value = 13.973799126637
$ws.Range("M19").Value = 52.186588921282
$ws.Range("N19").Value = -16.613418530351

# Row 20 - G.L.A.
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 128
$ws.Range("J20").Value = 125
$ws.Range("K20").Value = 2.4
$ws.Range("L20").Value = 47.126436781609
$ws.Range("M20").Value = -3.030303030303
$ws.Range("N20").Value = -91.208791208791

# Row 21 - TOTAL
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -27.027027027027
$ws.Range("G21").Value = 138
$ws.Range("H21").Value = -38.405797101449
$ws.Range("I21").Value = 1138
$ws.Range("J21").Value = 1167
$ws.Range("K21").Value = -2.485004284490
$ws.Range("L21").Value = 25.884955752212
$ws.Range("M21").Value = 17.927461139896
$ws.Range("N21").Value = -73.442240373395

# Row 22 - Transit
Set-TextValue "D22" "0" "C14"
Set-TextValue "E22" "***.*" "C14"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("M22").Value = -54.545454545454

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -28.947368421052
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 149
$ws.Range("H24").Value = -27.516778523489
$ws.Range("I24").Value = 1418
$ws.Range("J24").Value = 1505
$ws.Range("K24").Value = -5.780730897009
$ws.Range("L24").Value = 44.989775051124
$ws.Range("M24").Value = 74.201474201474

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 55.555555555555
$ws.Range("F25").Value = 58
$ws.Range("H25").Value = 48.717948717948
$ws.Range("I25").Value = 494
$ws.Range("J25").Value = 368
$ws.Range("K25").Value = 34.239130434782
$ws.Range("L25").Value = 42.774566473988
$ws.Range("M25").Value = 38.375350140056

# Row 26 - UCR Rape*
Set-TextValue "C26" "0" "C14"
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -66.666666666666
$ws.Range("J26").Value = 21
$ws.Range("K26").Value = -9.523809523809
$ws.Range("L26").Value = -9.523809523809

# Row 27 - Other Sex Crimes
$ws.Range("J27").Value = 42
$ws.Range("K27").Value = 7.142857142857
$ws.Range("L27").Value = -10

# Row 28 - Shooting Vic.
Set-TextValue "G28" "0" "C14"
Set-TextValue "H28" "***.*" "C14"

# Row 29 - Shooting Inc.
Set-TextValue "G29" "0" "C14"
Set-TextValue "H29" "***.*" "C14"
